$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.718.18"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.293.02"
$ws.Range("E3").Value = "  +5.47%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.60%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.292.30"
$ws.Range("E8").Value = "  +5.65%  "
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "3.838.62"
$ws.Range("E15").Value = "  +5.86%  "
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "3.294.66"
$ws.Range("E17").Value = "  +5.73%  "
$ws.Range("D18").Value = "63.782.51"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.96%  "
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.29%  "
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").Value = "0.0₃0749"
$ws.Range("E38").Value = "  +10.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "425.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "3.052.14"
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("E46").Value = "  +4.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("E51").Value = "  +2.55%  "
